# Re-format the single "questions = [...]" string (A2) as pretty-printed
# JSON-ish text, move it into A1 (dropping A1's old bold/bordered "0"
# placeholder value + style), and drop the now-empty second row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are a front-end developer working on a website for a client. You encounter a use-case where you need to render an element conditionally, hiding it under a certain condition. Which CSS rule should you apply to an element to ensure that it doesn't occupy any space on the page and is removed from the normal document flow?",
        "ques_type": 2,
        "options": [
            "display: none",
            "visibility: hidden",
            "opacity: 0",
            "position:absolute"
        ],
        "score": "display: none"
    },
    {
        "title": "You are a front-end developer working on a website for a client who is a professional photographer. You're creating a photo gallery using flexbox and want all thumbnails to be displayed in rows with equal height regardless of their original aspect ratio. Which of the following CSS properties should be applied to the thumbnails to achieve this effect?",
        "ques_type": 2,
        "options": [
            "align-items: stretch",
            "align-self: stretch",
            "flex: 1",
            "object-fit: cover"
        ],
        "score": "align-items: stretch"
    },
    {
        "title": "You are working on a webpage for a sports blog that has multiple article elements. Each article has an h2 heading followed by multiple h3 headings. You want to style only the first h3 heading that directly follows an h2 heading, without adding any classes or IDs. You decide to use an advanced CSS selector. Which selector will help you target the desired h3 element?",
        "ques_type": 2,
        "options": [
            "h2 + h3",
            "h2 ~ h3",
            "article h3:first-of-type",
            "article h3:nth-of-type(1)"
        ],
        "score": "h2 + h3"
    },
    {
        "title": "You are a web developer working on a company website. The company wants to add smooth transitions to the navigation menu items when users hover over them. You decide to use the transition-timing-function property in CSS to achieve this effect. Which of the following transition-timing-function values should you use to create a smooth and gradual transition effect?",
        "ques_type": 2,
        "options": [
            "Use ease-in to create a gradual transition effect with a slow start.",
            "Use cubic-bezier(0.5, 0, 0.5, 1) for precise transition control.",
            "Use linear to create a linear transition effect with a constant speed.",
            "Use step-start to create an immediate transition effect without any animation."
        ],
        "score": "Use cubic-bezier(0.5, 0, 0.5, 1) for precise transition control."
    }
]
'@

# Reset A1 to the default "Normal" style first (clears the bold font +
# thin border it had), then overwrite its value with the new text.
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = $newText

# The old question text lived in A2; remove that row entirely now that
# its content has moved up into A1.
$ws.Rows("2:2").Delete()

# Entering multi-line text auto-expands the row height; put it back to
# an auto (non-custom) height so row 1 keeps its default sizing.
$ws.Rows("1:1").AutoFit()
